$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same modified-date value for every data
# row (C2:C506). The value changes from 45175 (2023-09-06) to 45177
# (2023-09-08).
$ws.Range("C2:C506").Value = 45177
